# Update "想去人数" (want-to-go count) figures that were refreshed in the
# latest scrape run (gh-pages output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 13133   # 苏州·Good jump ACG中秋嘉年华动漫国潮文化节
$ws1.Range("F10").Value = 13088   # 苏州·I COME ACG动漫品牌博览会
$ws1.Range("F13").Value = 8784    # 苏州·理想乡动漫游戏展-两馆全开+三馆间通道
$ws1.Range("F14").Value = 7839    # 苏州·第四届-OCG国朝动漫游戏嘉年华
$ws1.Range("F23").Value = 191     # 苏州·第三届华盟国漫次元嘉年华

# Sheet "演出" (Performance)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 18       # 苏州·乐队番同人only live Band Set二次元乐队拼盘

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value  = 13133   # 苏州·Good jump ACG中秋嘉年华动漫国潮文化节
$ws4.Range("F11").Value = 13089   # 苏州·I COME ACG动漫品牌博览会
$ws4.Range("F14").Value = 8784    # 苏州·理想乡动漫游戏展-两馆全开+三馆间通道
$ws4.Range("F15").Value = 7839    # 苏州·第四届-OCG国朝动漫游戏嘉年华
$ws4.Range("F23").Value = 18      # 苏州·乐队番同人only live Band Set二次元乐队拼盘
$ws4.Range("F26").Value = 191     # 苏州·第三届华盟国漫次元嘉年华
